$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column C width from 23.71.. to 21.71.. (approx, closest achievable step)
$ws.Columns.Item(3).ColumnWidth = 20.83

# Row 2
$ws.Cells.Item(2, 1).Value = 300
$ws.Cells.Item(2, 2).Value = 76657
$ws.Cells.Item(2, 3).Value = 9210
$ws.Cells.Item(2, 8).Value = -55507
$ws.Cells.Item(2, 9).Value = -66207
$ws.Cells.Item(2, 10).Value = 253

# Row 3
$ws.Cells.Item(3, 1).Value = 46
$ws.Cells.Item(3, 2).Value = 2868
$ws.Cells.Item(3, 3).Value = 646
$ws.Cells.Item(3, 8).Value = -23774
$ws.Cells.Item(3, 9).Value = -869
$ws.Cells.Item(3, 10).Value = 371

# Row 4
$ws.Cells.Item(4, 1).Value = 126
$ws.Cells.Item(4, 2).Value = 6710
$ws.Cells.Item(4, 3).Value = 898
$ws.Cells.Item(4, 8).Value = 15501
$ws.Cells.Item(4, 9).Value = -7582
$ws.Cells.Item(4, 10).Value = 260

# Row 5
$ws.Cells.Item(5, 1).Value = 720
$ws.Cells.Item(5, 2).Value = 119922
$ws.Cells.Item(5, 3).Value = 9564
$ws.Cells.Item(5, 8).Value = 24118
$ws.Cells.Item(5, 9).Value = 79638
$ws.Cells.Item(5, 10).Value = 628

# Row 6
$ws.Cells.Item(6, 1).Value = 959
$ws.Cells.Item(6, 2).Value = 276787
$ws.Cells.Item(6, 3).Value = 34271
$ws.Cells.Item(6, 8).Value = 183236
$ws.Cells.Item(6, 9).Value = 343816
$ws.Cells.Item(6, 10).Value = 686

# Row 7
$ws.Cells.Item(7, 1).Value = 912
$ws.Cells.Item(7, 2).Value = 198551
$ws.Cells.Item(7, 3).Value = 58652
$ws.Cells.Item(7, 8).Value = 185909
$ws.Cells.Item(7, 9).Value = 346514
$ws.Cells.Item(7, 10).Value = 692

# Row 8
$ws.Cells.Item(8, 1).Value = 453
$ws.Cells.Item(8, 2).Value = 59529
$ws.Cells.Item(8, 3).Value = 10614
$ws.Cells.Item(8, 8).Value = 39813
$ws.Cells.Item(8, 9).Value = 92082
$ws.Cells.Item(8, 10).Value = 296

# Row 9
$ws.Cells.Item(9, 1).Value = 336
$ws.Cells.Item(9, 2).Value = 47835
$ws.Cells.Item(9, 3).Value = 26795
$ws.Cells.Item(9, 8).Value = 33834
$ws.Cells.Item(9, 9).Value = 71464
$ws.Cells.Item(9, 10).Value = 193

# Row 10
$ws.Cells.Item(10, 1).Value = 237
$ws.Cells.Item(10, 2).Value = 23477
$ws.Cells.Item(10, 3).Value = 5100
$ws.Cells.Item(10, 8).Value = -13358
$ws.Cells.Item(10, 9).Value = 13504
$ws.Cells.Item(10, 10).Value = 345

# Row 11
$ws.Cells.Item(11, 1).Value = 529
$ws.Cells.Item(11, 2).Value = 69644
$ws.Cells.Item(11, 3).Value = 8732
$ws.Cells.Item(11, 8).Value = 17371
$ws.Cells.Item(11, 9).Value = 58760
$ws.Cells.Item(11, 10).Value = 341

# Row 12
$ws.Cells.Item(12, 1).Value = 415
$ws.Cells.Item(12, 2).Value = 107979
$ws.Cells.Item(12, 3).Value = 61888
$ws.Cells.Item(12, 8).Value = -1914
$ws.Cells.Item(12, 9).Value = 32923
$ws.Cells.Item(12, 10).Value = 528

# Row 13
$ws.Cells.Item(13, 1).Value = 825
$ws.Cells.Item(13, 2).Value = 188100
$ws.Cells.Item(13, 3).Value = 93465
$ws.Cells.Item(13, 8).Value = -18612
$ws.Cells.Item(13, 9).Value = -331
$ws.Cells.Item(13, 10).Value = 543

# Row 14
$ws.Cells.Item(14, 1).Value = 1235
$ws.Cells.Item(14, 2).Value = 1317798
$ws.Cells.Item(14, 3).Value = 786198
$ws.Cells.Item(14, 8).Value = 523096
$ws.Cells.Item(14, 9).Value = 887427
$ws.Cells.Item(14, 10).Value = 842

# Row 15
$ws.Cells.Item(15, 1).Value = 517
$ws.Cells.Item(15, 2).Value = 200479
$ws.Cells.Item(15, 3).Value = 134647
$ws.Cells.Item(15, 8).Value = 184046
$ws.Cells.Item(15, 9).Value = 320207
$ws.Cells.Item(15, 10).Value = 769

# Row 16
$ws.Cells.Item(16, 1).Value = 94
$ws.Cells.Item(16, 2).Value = 6074
$ws.Cells.Item(16, 3).Value = 1775
$ws.Cells.Item(16, 8).Value = -35233
$ws.Cells.Item(16, 9).Value = -32265
$ws.Cells.Item(16, 10).Value = 455

# Row 17
$ws.Cells.Item(17, 1).Value = 122
$ws.Cells.Item(17, 2).Value = 10418
$ws.Cells.Item(17, 3).Value = 2338
$ws.Cells.Item(17, 8).Value = 9422
$ws.Cells.Item(17, 9).Value = 47883
$ws.Cells.Item(17, 10).Value = 343

# Row 18
$ws.Cells.Item(18, 1).Value = 244
$ws.Cells.Item(18, 2).Value = 28917
$ws.Cells.Item(18, 3).Value = 18388
$ws.Cells.Item(18, 8).Value = 54906
$ws.Cells.Item(18, 9).Value = 115841
$ws.Cells.Item(18, 10).Value = 690

# Row 19
$ws.Cells.Item(19, 1).Value = 383
$ws.Cells.Item(19, 2).Value = 99483
$ws.Cells.Item(19, 3).Value = 81724
$ws.Cells.Item(19, 8).Value = -19326
$ws.Cells.Item(19, 9).Value = 772
$ws.Cells.Item(19, 10).Value = 490

# Row 20
$ws.Cells.Item(20, 1).Value = 558
$ws.Cells.Item(20, 2).Value = 24722
$ws.Cells.Item(20, 3).Value = 12674
$ws.Cells.Item(20, 8).Value = 122322
$ws.Cells.Item(20, 9).Value = 208130
$ws.Cells.Item(20, 10).Value = 240

# Row 21
$ws.Cells.Item(21, 1).Value = 543
$ws.Cells.Item(21, 2).Value = 31199
$ws.Cells.Item(21, 3).Value = 9198
$ws.Cells.Item(21, 8).Value = -8380
$ws.Cells.Item(21, 9).Value = 26871
$ws.Cells.Item(21, 10).Value = 457

# Row 22
$ws.Cells.Item(22, 1).Value = 633
$ws.Cells.Item(22, 2).Value = 227985
$ws.Cells.Item(22, 3).Value = 185155
$ws.Cells.Item(22, 8).Value = 334463
$ws.Cells.Item(22, 9).Value = 586557
$ws.Cells.Item(22, 10).Value = 765
